$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sss"
$ws.Range("B2").Value = "sss"
$ws.Range("B3").Value = "sss"
$ws.Range("C3").Value = "ss"

$ws.Range("C3").Select()
